# Update nombre_aides (column C) and montant_total (column D) for the 2020-07-29 data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, new value for column C (nombre_aides), new value for column D (montant_total)
$updates = @(
    ,(2, 36942, 53465825)
    ,(3, 89416, 131136455)
    ,(4, 30653, 45405596)
    ,(5, 8535, 12687076)
    ,(6, 1917, 2849093)
    ,(7, 143, 209593)
    ,(11, 40346, 54788820)
    ,(12, 9436, 13651460)
    ,(13, 25554, 37486616)
    ,(14, 8203, 12175041)
    ,(15, 2115, 3144883)
    ,(16, 402, 592123)
    ,(19, 9999, 13267360)
    ,(20, 13187, 19049037)
    ,(21, 31241, 45862472)
    ,(22, 10107, 15028997)
    ,(23, 2586, 3849263)
    ,(24, 487, 724592)
    ,(26, 11466, 15340473)
    ,(27, 7501, 10869851)
    ,(28, 22163, 32531869)
    ,(29, 7690, 11444292)
    ,(30, 1926, 2873500)
    ,(31, 350, 522415)
    ,(32, 28, 41893)
    ,(33, 8165, 10792385)
    ,(34, 3141, 4533245)
    ,(35, 7618, 11126788)
    ,(36, 3097, 4589788)
    ,(37, 807, 1203263)
    ,(38, 146, 217232)
    ,(40, 2378, 3211861)
    ,(41, 16926, 24485970)
    ,(42, 50265, 73711701)
    ,(43, 18731, 27824821)
    ,(44, 5516, 8213978)
    ,(45, 1143, 1705292)
    ,(46, 59, 86848)
    ,(49, 16388, 21840938)
    ,(50, 1904, 2762335)
    ,(51, 6571, 9666330)
    ,(52, 2264, 3381324)
    ,(53, 739, 1103805)
    ,(56, 6373, 8782602)
    ,(57, 859, 1259834)
    ,(58, 2134, 3167037)
    ,(59, 876, 1303501)
    ,(60, 296, 443758)
    ,(61, 93, 139500)
    ,(63, 1261, 1786073)
    ,(64, 15098, 21815864)
    ,(65, 44065, 64501850)
    ,(66, 15488, 23020061)
    ,(67, 4490, 6688675)
    ,(68, 900, 1339596)
    ,(72, 14832, 19575002)
    ,(73, 49847, 72563523)
    ,(74, 142334, 209740818)
    ,(75, 62171, 92660048)
    ,(76, 19805, 29592046)
    ,(77, 4608, 6884906)
    ,(78, 245, 362670)
    ,(84, 49397, 67344151)
    ,(85, 4477, 6488386)
    ,(86, 11289, 16587040)
    ,(87, 3814, 5684165)
    ,(92, 5230, 7045662)
    ,(93, 1526, 2199539)
    ,(94, 4990, 7351737)
    ,(95, 1893, 2820999)
    ,(96, 666, 997960)
    ,(97, 173, 258613)
    ,(100, 3395, 4501133)
    ,(101, 578, 862464)
    ,(102, 337, 503130)
    ,(103, 117, 175500)
    ,(106, 10584, 15364107)
    ,(107, 28818, 42344191)
    ,(108, 9640, 14336657)
    ,(109, 2651, 3952707)
    ,(110, 473, 705796)
    ,(111, 47, 70500)
    ,(113, 9619, 12722009)
    ,(114, 29858, 43080172)
    ,(115, 65156, 95379259)
    ,(116, 21049, 31288222)
    ,(117, 5927, 8831780)
    ,(118, 1094, 1635506)
    ,(119, 74, 108420)
    ,(122, 25282, 33801580)
    ,(123, 35085, 50661613)
    ,(124, 75315, 110186147)
    ,(125, 23436, 34787405)
    ,(126, 6253, 9292307)
    ,(127, 1180, 1754411)
    ,(131, 30955, 41142163)
    ,(132, 13021, 18853395)
    ,(133, 31874, 46830012)
    ,(134, 11336, 16845362)
    ,(135, 2903, 4329081)
    ,(137, 33, 48325)
    ,(139, 10640, 14212451)
    ,(140, 34349, 49626209)
    ,(141, 79826, 116980124)
    ,(142, 23925, 35551412)
    ,(143, 6257, 9338726)
    ,(144, 1391, 2069692)
    ,(147, 28578, 38617305)
)

foreach ($u in $updates) {
    $row = $u[0]
    $newC = $u[1]
    $newD = $u[2]
    $ws.Cells.Item($row, 3).Value = $newC
    $ws.Cells.Item($row, 4).Value = $newD
}
